$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.1477750351608889
$ws.Range("J2").Value = 0.1477750351608889
$ws.Range("M2").Value = 0.794582
$ws.Range("N2").Value = 2.383746
$ws.Range("O2").Value = 0.03449752952410986
$ws.Range("P2").Value = 0.03449752952410985
$ws.Range("Q2").Value = 0.057940124858
$ws.Range("R2").Value = 0.521461123722
$ws.Range("S2").Value = 0.005097873638389138
$ws.Range("T2").Value = 0.005097873638389137

$ws.Range("I3").Value = 0.1477750351608889
$ws.Range("J3").Value = 0.1477750351608889
$ws.Range("O3").Value = 0.8945489325574519
$ws.Range("P3").Value = 0.8945489325574517
$ws.Range("S3").Value = 0.1321919999618131
$ws.Range("T3").Value = 0.1321919999618131

$ws.Range("I4").Value = 0.1477750351608889
$ws.Range("J4").Value = 0.1477750351608889
$ws.Range("M4").Value = 0.2871986666666667
$ws.Range("N4").Value = 0.8615959999999999
$ws.Range("O4").Value = 0.01246900191876775
$ws.Range("P4").Value = 0.01246900191876775
$ws.Range("Q4").Value = 0.02094223957466667
$ws.Range("R4").Value = 0.188480156172
$ws.Range("S4").Value = 0.001842607196967096
$ws.Range("T4").Value = 0.001842607196967096

$ws.Range("I5").Value = 0.1477750351608889
$ws.Range("J5").Value = 0.1477750351608889
$ws.Range("M5").Value = 1.149534666666667
$ws.Range("N5").Value = 3.448604
$ws.Range("O5").Value = 0.04990813547540859
$ws.Range("P5").Value = 0.04990813547540859
$ws.Range("Q5").Value = 0.08382291835866666
$ws.Range("R5").Value = 0.754406265228
$ws.Range("S5").Value = 0.007375176474692912
$ws.Range("T5").Value = 0.007375176474692912

$ws.Range("I6").Value = 0.1477750351608889
$ws.Range("J6").Value = 0.1477750351608889
$ws.Range("M6").Value = 0.1975403333333333
$ws.Range("N6").Value = 0.5926210000000001
$ws.Range("O6").Value = 0.008576400524262026
$ws.Range("P6").Value = 0.008576400524262026
$ws.Range("Q6").Value = 0.01440444356633333
$ws.Range("R6").Value = 0.129639992097
$ws.Range("S6").Value = 0.001267377889026687
$ws.Range("T6").Value = 0.001267377889026687

$ws.Range("G7").Value = 0.420527
$ws.Range("H7").Value = 1.261581
$ws.Range("I7").Value = 0.852224964839111
$ws.Range("J7").Value = 0.852224964839111
$ws.Range("M7").Value = 0.794582
$ws.Range("N7").Value = 2.383746
$ws.Range("O7").Value = 0.03449752952410986
$ws.Range("P7").Value = 0.03449752952410985
$ws.Range("Q7").Value = 0.3341431847140001
$ws.Range("R7").Value = 3.007288662426
$ws.Range("S7").Value = 0.02939965588572072
$ws.Range("T7").Value = 0.02939965588572071

$ws.Range("G8").Value = 0.420527
$ws.Range("H8").Value = 1.261581
$ws.Range("I8").Value = 0.852224964839111
$ws.Range("J8").Value = 0.852224964839111
$ws.Range("O8").Value = 0.8945489325574519
$ws.Range("P8").Value = 0.8945489325574517
$ws.Range("Q8").Value = 8.664603910212001
$ws.Range("R8").Value = 77.98143519190801
$ws.Range("S8").Value = 0.7623569325956387
$ws.Range("T8").Value = 0.7623569325956386

$ws.Range("G9").Value = 0.420527
$ws.Range("H9").Value = 1.261581
$ws.Range("I9").Value = 0.852224964839111
$ws.Range("J9").Value = 0.852224964839111
$ws.Range("M9").Value = 0.2871986666666667
$ws.Range("N9").Value = 0.8615959999999999
$ws.Range("O9").Value = 0.01246900191876775
$ws.Range("P9").Value = 0.01246900191876775
$ws.Range("Q9").Value = 0.1207747936973333
$ws.Range("R9").Value = 1.086973143276
$ws.Range("S9").Value = 0.01062639472180066
$ws.Range("T9").Value = 0.01062639472180065

$ws.Range("G10").Value = 0.420527
$ws.Range("H10").Value = 1.261581
$ws.Range("I10").Value = 0.852224964839111
$ws.Range("J10").Value = 0.852224964839111
$ws.Range("M10").Value = 1.149534666666667
$ws.Range("N10").Value = 3.448604
$ws.Range("O10").Value = 0.04990813547540859
$ws.Range("P10").Value = 0.04990813547540859
$ws.Range("Q10").Value = 0.4834103647693334
$ws.Range("R10").Value = 4.350693282924
$ws.Range("S10").Value = 0.04253295900071568
$ws.Range("T10").Value = 0.04253295900071568

$ws.Range("G11").Value = 0.420527
$ws.Range("H11").Value = 1.261581
$ws.Range("I11").Value = 0.852224964839111
$ws.Range("J11").Value = 0.852224964839111
$ws.Range("M11").Value = 0.1975403333333333
$ws.Range("N11").Value = 0.5926210000000001
$ws.Range("O11").Value = 0.008576400524262026
$ws.Range("P11").Value = 0.008576400524262026
$ws.Range("Q11").Value = 0.08307104375566668
$ws.Range("R11").Value = 0.7476393938010001
$ws.Range("S11").Value = 0.007309022635235338
$ws.Range("T11").Value = 0.007309022635235338
